{"js": "// Resume content update:\n//   1. Bold the tool names \"Python\", \"Microsoft Power BI\", \"Tableau\" in the\n//      \"Data Exploration, Analysis, and Visualization\" bullet.\n//   2. Move the (hidden) \"_GoBack\" bookmark from the end of the\n//      \"Insurance, Risk Management, Financial:\" heading up to right after\n//      \"...system performance and \" in the \"How I Add Value\" bullet.\n\nconst body = context.document.body;\n\n// --- 1. Bold the tool names -------------------------------------------------\nconst toolNames = [\"Python\", \"Microsoft Power BI\", \"Tableau\"];\nfor (const name of toolNames) {\n  const found = body.search(name, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length > 0) {\n    found.items[0].font.set({ bold: true });\n  }\n}\nawait context.sync();\n\n// --- 2. Relocate the \"_GoBack\" bookmark ------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst anchor = body.search(\"system performance and \", { matchCase: true });\nanchor.load(\"items\");\nawait context.sync();\n\nif (anchor.items.length > 0) {\n  const insertionPoint = anchor.items[0].getRange(\"After\");\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Resume content update:\n#   1. Bold the tool names \"Python\", \"Microsoft Power BI\", \"Tableau\" in the\n#      \"Data Exploration, Analysis, and Visualization\" bullet.\n#   2. Move the (hidden) \"_GoBack\" bookmark from the end of the\n#      \"Insurance, Risk Management, Financial:\" heading up to right after\n#      \"...system performance and \" in the \"How I Add Value\" bullet.\n\n$d = $word.ActiveDocument\n\n# --- 1. Bold the tool names -------------------------------------------------\n$toolNames = @(\"Python\", \"Microsoft Power BI\", \"Tableau\")\nforeach ($name in $toolNames) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $name\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $found = $rng.Find.Execute()\n    if ($found) {\n        $rng.Bold = 1\n    }\n}\n\n# --- 2. Relocate the \"_GoBack\" bookmark ------------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$anchor = $d.Content\n$anchor.Find.ClearFormatting()\n$anchor.Find.Text = \"system performance and \"\n$anchor.Find.MatchCase = $true\n$anchor.Find.MatchWholeWord = $false\n$found = $anchor.Find.Execute()\nif ($found) {\n    $insertPoint = $d.Range($anchor.End, $anchor.End)\n    $d.Bookmarks.Add(\"_GoBack\", $insertPoint)\n}\n"}
